# Generate Report for Handback
# Updates the handback-status report with a freshly generated handback
# entry: the old in-flight GUID (e52db018-...) resolved/renamed to
# 58dc30ff-4c75-47e5-95bf-406b1a6b723e, and a brand-new second entry
# (5ec2340e-2472-47a6-8ec5-01d3c83f9145) was appended as row 3 on every
# sheet.

$wb = $excel.ActiveWorkbook

$oldGuid  = "e52db018-6830-4abf-ba72-dd23dfc01521"
$newGuid1 = "58dc30ff-4c75-47e5-95bf-406b1a6b723e"
$newGuid2 = "5ec2340e-2472-47a6-8ec5-01d3c83f9145"

$oldHash  = "c7f9f1f55fc52c0b745fea7ef0e3f7206b93ab5a"
$newHash1 = "f5be56a1146a15bd60a4a7cfdf825f8351aff8e5"
$newHash2 = "8a631ac5f48729918fecb1585a2f74ac01b3fb22"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2: the in-flight GUID is now resolved -> rename file + bump date.
$ov.Range("A2").Value = "'" + $newGuid1 + ".md"
$ov.Range("B2").Value = "'e2e\" + $newGuid1 + ".md"
$ov.Range("G2").Value = "'2016-08-21 03:04:30"

# The existing hyperlink still points at the old filename -> replace it.
$ov.Range("B2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/" + $newGuid1 + ".md", "", "", "e2e\" + $newGuid1 + ".md") | Out-Null

# Row 3: brand-new handback entry.
$ov.Range("A3").Value = "'" + $newGuid2 + ".md"
$ov.Range("B3").Value = "'e2e\" + $newGuid2 + ".md"
$ov.Range("C3").Value = "'.md"
$ov.Range("E3").Value = "'Handed back: in sync with en-US"
$ov.Range("F3").Value = "'Handed back: in sync with en-US"
$ov.Range("G3").Value = "'2016-08-21 03:04:30"

$ov.Range("B3").NumberFormat = $ov.Range("B2").NumberFormat
$ov.Range("B3").Font.Underline = $ov.Range("B2").Font.Underline
$ov.Range("B3").Font.Color = $ov.Range("B2").Font.Color
$ov.Range("G3").NumberFormat = $ov.Range("G2").NumberFormat

$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/" + $newGuid2 + ".md", "", "", "e2e\" + $newGuid2 + ".md") | Out-Null

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de" share the same 16-column layout; only the
# xlf suffix (zh-cn.xlf / de-de.xlf) and the "Correspond Handoff
# Datetime" / "Correspond Handback DateTime" stamps differ between them.
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Suffix = "zh-cn"; HandoffDate2 = "2016-08-21 03:04:25"; HandbackDate2 = "2016-08-21 03:04:41"; HandoffDate3 = "2016-08-21 03:04:25"; HandbackDate3 = "2016-08-21 03:04:41"; UrlOrg = "ol-test0-zhcn"; UrlCommit = "ec033b5b126227c736920f0dac3292d84b1e655f" },
    @{ Name = "de-de"; Suffix = "de-de"; HandoffDate2 = "2016-08-21 03:04:30"; HandbackDate2 = "2016-08-21 03:04:48"; HandoffDate3 = "2016-08-21 03:04:30"; HandbackDate3 = "2016-08-21 03:04:48"; UrlOrg = "ol-test0-dede"; UrlCommit = "0f1da85f5f48929c8121f49c6701fbac71fff4cf" }
)

foreach ($info in $langSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 2: same physical file, now carries the resolved GUID + new hash/dates.
    $ws.Range("A2").Value = "'" + $newGuid1 + ".md"
    $ws.Range("G2").Value = "'" + $newGuid1 + "." + $newHash1 + "." + $info.Suffix + ".xlf"
    $ws.Range("H2").Value = "'" + $info.HandoffDate2
    $ws.Range("I2").Value = "'" + $newGuid1 + ".md"
    $ws.Range("J2").Value = "'" + $newGuid1 + "." + $newHash1 + "." + $info.Suffix + ".xlf"
    $ws.Range("K2").Value = "'" + $info.HandbackDate2

    # The existing hyperlinks still point at the old filename -> replace them.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("I2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/" + $newGuid1 + ".md", "", "", $newGuid1 + ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/" + $info.UrlOrg + "/blob/" + $info.UrlCommit + "/e2e/" + $newGuid1 + ".md", "", "", $newGuid1 + ".md") | Out-Null

    # Row 3: new handback entry, mirrors row 2's layout/formatting.
    $ws.Range("A3").Value = "'" + $newGuid2 + ".md"
    $ws.Range("B3").Value = "'.md"
    $ws.Range("C3").Value = "'Handed back: in sync with en-US"
    $ws.Range("D3").Value = "'e2e"
    $ws.Range("E3").Value = "'ht"
    $ws.Range("F3").Value = "'True"
    $ws.Range("G3").Value = "'" + $newGuid2 + "." + $newHash2 + "." + $info.Suffix + ".xlf"
    $ws.Range("H3").Value = "'" + $info.HandoffDate3
    $ws.Range("I3").Value = "'" + $newGuid2 + ".md"
    $ws.Range("J3").Value = "'" + $newGuid2 + "." + $newHash2 + "." + $info.Suffix + ".xlf"
    $ws.Range("K3").Value = "'" + $info.HandbackDate3
    $ws.Range("L3").Value = "'"
    $ws.Range("M3").Value = "'True"
    $ws.Range("N3").Value = "'"
    $ws.Range("O3").Value = "'False"
    $ws.Range("P3").Value = "'"

    # Copy row 2's formatting onto row 3 so dates/hyperlinks render the same way.
    $ws.Range("A2:P2").Copy()
    $ws.Range("A3:P3").PasteSpecial(-4122)
    $ws.Range("A3").Value = "'" + $newGuid2 + ".md"
    $ws.Range("B3").Value = "'.md"
    $ws.Range("C3").Value = "'Handed back: in sync with en-US"
    $ws.Range("D3").Value = "'e2e"
    $ws.Range("E3").Value = "'ht"
    $ws.Range("F3").Value = "'True"
    $ws.Range("G3").Value = "'" + $newGuid2 + "." + $newHash2 + "." + $info.Suffix + ".xlf"
    $ws.Range("H3").Value = "'" + $info.HandoffDate3
    $ws.Range("I3").Value = "'" + $newGuid2 + ".md"
    $ws.Range("J3").Value = "'" + $newGuid2 + "." + $newHash2 + "." + $info.Suffix + ".xlf"
    $ws.Range("K3").Value = "'" + $info.HandbackDate3
    $ws.Range("L3").Value = "'"
    $ws.Range("M3").Value = "'True"
    $ws.Range("N3").Value = "'"
    $ws.Range("O3").Value = "'False"
    $ws.Range("P3").Value = "'"

    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1270c9436dbc1cf8a2c4ac200d257e5a5cb68c/e2e/" + $newGuid2 + ".md", "", "", $newGuid2 + ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/" + $info.UrlOrg + "/blob/" + $info.UrlCommit + "/e2e/" + $newGuid2 + ".md", "", "", $newGuid2 + ".md") | Out-Null

    $t = $ws.ListObjects.Item(1)
    $t.Resize($ws.Range("A1:P3"))
}
